$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "addTaskTest"

$ws.Range("A1").Value = "task"
$ws.Range("B1").Value = "description"
$ws.Range("A2").Value = "Null Tasks"
$ws.Range("B2").Value = "Automation skills"
$ws.Range("A3").Value = "2nd task"
$ws.Range("B3").Value = "null pointer"

$ws.Range("C6").Select()
